# Update "想去人数" (F column) values across all 4 sheets to reflect new
# scrape snapshot (gh-pages output generated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2363
$ws.Range("F3").Value = 543
$ws.Range("F5").Value = 356
$ws.Range("F6").Value = 356
$ws.Range("F7").Value = 584
$ws.Range("F9").Value = 798
$ws.Range("F11").Value = 823
$ws.Range("F12").Value = 390
$ws.Range("F13").Value = 101
$ws.Range("F14").Value = 400
$ws.Range("F16").Value = 1032
$ws.Range("F17").Value = 21316
$ws.Range("F18").Value = 850
$ws.Range("F19").Value = 81
$ws.Range("F20").Value = 263
$ws.Range("F23").Value = 169
$ws.Range("F25").Value = 17
$ws.Range("F26").Value = 237
$ws.Range("F28").Value = 360
$ws.Range("F29").Value = 159

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 185
$ws.Range("F5").Value = 89
$ws.Range("F8").Value = 3441
$ws.Range("F10").Value = 105
$ws.Range("F16").Value = 3895

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 114
$ws.Range("F4").Value = 629
$ws.Range("F5").Value = 214

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 114
$ws.Range("F5").Value = 2363
$ws.Range("F6").Value = 629
$ws.Range("F7").Value = 543
$ws.Range("F9").Value = 356
$ws.Range("F10").Value = 356
$ws.Range("F11").Value = 584
$ws.Range("F12").Value = 185
$ws.Range("F15").Value = 89
$ws.Range("F17").Value = 214
$ws.Range("F18").Value = 798
$ws.Range("F20").Value = 823
$ws.Range("F21").Value = 390
$ws.Range("F22").Value = 101
$ws.Range("F23").Value = 400
$ws.Range("F25").Value = 1032
$ws.Range("F26").Value = 21316
$ws.Range("F28").Value = 3441
$ws.Range("F30").Value = 105
$ws.Range("F32").Value = 850
$ws.Range("F33").Value = 81
$ws.Range("F34").Value = 263
$ws.Range("F39").Value = 169
$ws.Range("F41").Value = 17
$ws.Range("F44").Value = 237
$ws.Range("F46").Value = 360
$ws.Range("F47").Value = 159
$ws.Range("F48").Value = 3895

